# edit.ps1 - Update '想去人数' (F column, number of people interested) values
# across all 4 worksheets (展览, 演出, 本地生活, 全部类型), matching a refreshed
# data pull (gh-pages output regenerated at commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Auto-generated edit script: update '想去人数' (F column) values across sheets 1, 2, 3, 4

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 6399
$ws.Cells.Item(4, 6).Value = 1026
$ws.Cells.Item(5, 6).Value = 627
$ws.Cells.Item(6, 6).Value = 1409
$ws.Cells.Item(8, 6).Value = 3
$ws.Cells.Item(9, 6).Value = 506
$ws.Cells.Item(10, 6).Value = 2049
$ws.Cells.Item(11, 6).Value = 442
$ws.Cells.Item(13, 6).Value = 214
$ws.Cells.Item(14, 6).Value = 103
$ws.Cells.Item(15, 6).Value = 216
$ws.Cells.Item(16, 6).Value = 1024
$ws.Cells.Item(17, 6).Value = 388
$ws.Cells.Item(19, 6).Value = 141
$ws.Cells.Item(20, 6).Value = 3910
$ws.Cells.Item(21, 6).Value = 1205
$ws.Cells.Item(22, 6).Value = 3093
$ws.Cells.Item(23, 6).Value = 300
$ws.Cells.Item(24, 6).Value = 66
$ws.Cells.Item(25, 6).Value = 2642
$ws.Cells.Item(26, 6).Value = 2642
$ws.Cells.Item(27, 6).Value = 4465
$ws.Cells.Item(29, 6).Value = 946
$ws.Cells.Item(30, 6).Value = 500
$ws.Cells.Item(31, 6).Value = 2971
$ws.Cells.Item(32, 6).Value = 264
$ws.Cells.Item(33, 6).Value = 26
$ws.Cells.Item(34, 6).Value = 74
$ws.Cells.Item(35, 6).Value = 56
$ws.Cells.Item(36, 6).Value = 42
$ws.Cells.Item(37, 6).Value = 1078
$ws.Cells.Item(38, 6).Value = 1327
$ws.Cells.Item(39, 6).Value = 91
$ws.Cells.Item(40, 6).Value = 1177
$ws.Cells.Item(41, 6).Value = 764
$ws.Cells.Item(42, 6).Value = 3
$ws.Cells.Item(43, 6).Value = 700
$ws.Cells.Item(44, 6).Value = 460
$ws.Cells.Item(45, 6).Value = 37
$ws.Cells.Item(46, 6).Value = 163
$ws.Cells.Item(47, 6).Value = 19
$ws.Cells.Item(48, 6).Value = 62
$ws.Cells.Item(49, 6).Value = 335
$ws.Cells.Item(50, 6).Value = 3645

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 26
$ws.Cells.Item(10, 6).Value = 937
$ws.Cells.Item(13, 6).Value = 3
$ws.Cells.Item(18, 6).Value = 4
$ws.Cells.Item(19, 6).Value = 4

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 954

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 6399
$ws.Cells.Item(3, 6).Value = 26
$ws.Cells.Item(5, 6).Value = 627
$ws.Cells.Item(6, 6).Value = 1409
$ws.Cells.Item(8, 6).Value = 3
$ws.Cells.Item(9, 6).Value = 506
$ws.Cells.Item(11, 6).Value = 2049
$ws.Cells.Item(12, 6).Value = 442
$ws.Cells.Item(14, 6).Value = 214
$ws.Cells.Item(15, 6).Value = 937
$ws.Cells.Item(17, 6).Value = 103
$ws.Cells.Item(18, 6).Value = 216
$ws.Cells.Item(19, 6).Value = 1024
$ws.Cells.Item(21, 6).Value = 388
$ws.Cells.Item(22, 6).Value = 141
$ws.Cells.Item(23, 6).Value = 3910
$ws.Cells.Item(25, 6).Value = 1205
$ws.Cells.Item(26, 6).Value = 4
$ws.Cells.Item(27, 6).Value = 3094
$ws.Cells.Item(28, 6).Value = 2642
$ws.Cells.Item(29, 6).Value = 2642
$ws.Cells.Item(30, 6).Value = 4466
$ws.Cells.Item(31, 6).Value = 946
$ws.Cells.Item(32, 6).Value = 2971
$ws.Cells.Item(33, 6).Value = 264
$ws.Cells.Item(34, 6).Value = 42
$ws.Cells.Item(35, 6).Value = 1078
$ws.Cells.Item(36, 6).Value = 1327
$ws.Cells.Item(37, 6).Value = 91
$ws.Cells.Item(38, 6).Value = 1177
$ws.Cells.Item(39, 6).Value = 764
$ws.Cells.Item(41, 6).Value = 460
$ws.Cells.Item(44, 6).Value = 37
$ws.Cells.Item(46, 6).Value = 163
$ws.Cells.Item(47, 6).Value = 19
$ws.Cells.Item(48, 6).Value = 62
$ws.Cells.Item(49, 6).Value = 335
$ws.Cells.Item(50, 6).Value = 3645
